$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 19
$ws.Range("H2").Value = 82.61

# Row 3
$ws.Range("F3").Value = 379
$ws.Range("H3").Value = 96.19

# Row 6
$ws.Range("F6").Value = 46
$ws.Range("H6").Value = 92

# Row 7
$ws.Range("F7").Value = 60
$ws.Range("H7").Value = 96.77

# Row 8
$ws.Range("F8").Value = 89
$ws.Range("H8").Value = 95.7

# Row 10
$ws.Range("F10").Value = 75
$ws.Range("H10").Value = 89.29

# Row 11
$ws.Range("F11").Value = 79
$ws.Range("H11").Value = 90.8

# Row 12
$ws.Range("F12").Value = 100
$ws.Range("G12").Value = 137
$ws.Range("H12").Value = 72.99

# Row 13
$ws.Range("F13").Value = 107
$ws.Range("G13").Value = 127
$ws.Range("H13").Value = 84.25

# Row 14
$ws.Range("F14").Value = 174
$ws.Range("G14").Value = 206
$ws.Range("H14").Value = 84.47

# Row 15
$ws.Range("F15").Value = 127
$ws.Range("H15").Value = 100

# Row 16
$ws.Range("F16").Value = 99
$ws.Range("G16").Value = 118
$ws.Range("H16").Value = 83.9

# Row 19
$ws.Range("F19").Value = 25
$ws.Range("H19").Value = 65.79

# Row 20
$ws.Range("F20").Value = 48
$ws.Range("H20").Value = 97.96

# Row 21
$ws.Range("F21").Value = 126
$ws.Range("H21").Value = 93.33

# Row 22
$ws.Range("F22").Value = 46
$ws.Range("H22").Value = 80.7

# Row 23
$ws.Range("F23").Value = 416
$ws.Range("G23").Value = 491
$ws.Range("H23").Value = 84.73

# Row 24
$ws.Range("F24").Value = 39
$ws.Range("H24").Value = 73.58

# Row 25
$ws.Range("F25").Value = 28
$ws.Range("H25").Value = 96.55

# Row 26
$ws.Range("F26").Value = 76
$ws.Range("H26").Value = 76.77

# Row 27
$ws.Range("F27").Value = 29
$ws.Range("H27").Value = 87.88

# Row 29
$ws.Range("F29").Value = 10
$ws.Range("H29").Value = 83.33

# Row 30
$ws.Range("F30").Value = 166
$ws.Range("G30").Value = 241
$ws.Range("H30").Value = 68.88

# Row 31
$ws.Range("F31").Value = 101
$ws.Range("G31").Value = 143
$ws.Range("H31").Value = 70.63

# Row 32
$ws.Range("F32").Value = 165
$ws.Range("H32").Value = 87.77

# Row 33
$ws.Range("F33").Value = 114
$ws.Range("H33").Value = 90.48

# Row 34
$ws.Range("F34").Value = 83
$ws.Range("H34").Value = 84.69

# Row 35
$ws.Range("F35").Value = 23
$ws.Range("H35").Value = 67.65

# Row 36
$ws.Range("F36").Value = 40
$ws.Range("H36").Value = 90.91

# Row 37
$ws.Range("F37").Value = 70
$ws.Range("G37").Value = 78
$ws.Range("H37").Value = 89.74

# Row 38
$ws.Range("F38").Value = 121
$ws.Range("H38").Value = 80.67

# Row 39
$ws.Range("F39").Value = 74
$ws.Range("H39").Value = 79.57

# Row 40
$ws.Range("F40").Value = 17
$ws.Range("H40").Value = 65.38

# Row 41
$ws.Range("F41").Value = 62
$ws.Range("H41").Value = 80.52

# Row 42
$ws.Range("F42").Value = 125
$ws.Range("H42").Value = 71.84

# Row 43
$ws.Range("F43").Value = 97
$ws.Range("H43").Value = 83.62

# Row 44
$ws.Range("F44").Value = 171
$ws.Range("G44").Value = 202
$ws.Range("H44").Value = 84.65

# Row 45
$ws.Range("F45").Value = 119
$ws.Range("H45").Value = 57.77

# Row 46
$ws.Range("F46").Value = 8
$ws.Range("H46").Value = 72.73

# Row 47
$ws.Range("F47").Value = 58
$ws.Range("H47").Value = 85.29

# Row 48
$ws.Range("F48").Value = 90
$ws.Range("G48").Value = 104
$ws.Range("H48").Value = 86.54

# Row 49
$ws.Range("F49").Value = 10
$ws.Range("H49").Value = 71.43

# Row 50
$ws.Range("F50").Value = 69
$ws.Range("H50").Value = 80.23

# Row 51
$ws.Range("G51").Value = 56
$ws.Range("H51").Value = 82.14

# Row 52
$ws.Range("F52").Value = 104
$ws.Range("H52").Value = 79.39

# Row 54
$ws.Range("F54").Value = 2
$ws.Range("H54").Value = 66.67

# Row 57
$ws.Range("F57").Value = 64
$ws.Range("G57").Value = 65
$ws.Range("H57").Value = 98.46

# Row 58
$ws.Range("F58").Value = 2
$ws.Range("H58").Value = 100
